$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Print area: shrink from $A$1:$C$6 to $A$1:$C$5
$ws.PageSetup.PrintArea = '$A$1:$C$5'

# 2. Update shared string text: "Manager, Mike" -> "Manager, Mike \n(PMP)"
$ws.Range("A2").Value = "Manager, Mike " + [char]10 + "(PMP)"

# 3. Move the SUM formula from C5 to A5 and merge A5:C5
#    (do this before re-formatting C2:C4 so A5 doesn't inherit the
#    percentage number format from its precedent cells)
$ws.Range("C5").ClearContents()
$ws.Range("A5").Formula = "=SUM(C2:C4)"
$ws.Range("A5:C5").Merge()

# 4. Wrap the resource names in column A (rows 2-4) and give those rows
#    extra height to fit the two-line "Manager, Mike (PMP)" text
foreach ($addr in @("A2", "A3", "A4")) {
    $ws.Range($addr).WrapText = $true
}
$ws.Rows("2:4").RowHeight = 30

# 5. Convert the Utilization% values to fractions with a 0.0% number format
#    and left-align them
$ws.Range("C2").Value = 0.76
$ws.Range("C3").Value = 0.95
$ws.Range("C4").Value = 0.88
foreach ($addr in @("C2", "C3", "C4")) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "0.0%"
    $cell.HorizontalAlignment = -4131
}

# 6. Auto-size the columns to fit their new content
$ws.Columns("A:C").AutoFit()
